$d = $word.ActiveDocument

$replacements = @(
    @("37÷6=6, 1", "19÷9=2, 1"),
    @("56÷6=9, 2", "58÷6=9, 4"),
    @("12÷7=1, 5", "60÷2=30, 0"),
    @("89÷8=11, 1", "77÷5=15, 2"),
    @("99÷8=12, 3", "10÷4=2, 2"),
    @("81÷2=40, 1", "32÷8=4, 0"),
    @("26÷4=6, 2", "94÷2=47, 0"),
    @("56÷3=18, 2", "96÷6=16, 0"),
    @("27÷2=13, 1", "37÷3=12, 1"),
    @("96÷3=32, 0", "82÷6=13, 4"),
    @("69÷4=17, 1", "55÷2=27, 1"),
    @("89÷7=12, 5", "77÷2=38, 1"),
    @("12÷6=2, 0", "54÷6=9, 0"),
    @("44÷8=5, 4", "76÷6=12, 4"),
    @("17÷7=2, 3", "28÷9=3, 1"),
    @("74÷4=18, 2", "97÷3=32, 1"),
    @("59÷4=14, 3", "56÷8=7, 0"),
    @("93÷2=46, 1", "81÷4=20, 1"),
    @("17÷5=3, 2", "70÷9=7, 7"),
    @("40÷4=10, 0", "13÷3=4, 1"),
    @("46÷6=7, 4", "53÷2=26, 1"),
    @("68÷5=13, 3", "54÷9=6, 0"),
    @("70÷5=14, 0", "82÷2=41, 0"),
    @("54÷8=6, 6", "75÷9=8, 3"),
    @("63÷3=21, 0", "23÷8=2, 7")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
